$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.878.93'
$ws.Range("E2").Value = '  +5.17%  '

# Row 3
$ws.Range("D3").Value = '2.261.76'
$ws.Range("E3").Value = '  +2.48%  '

# Row 4
$ws.Range("E4").Value = '  +0.30%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.632'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.58%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.21'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.48%  '

# Row 8
$ws.Range("E8").Value = '  +0.13%  '

# Row 9
$ws.Range("E9").Value = '  +7.14%  '

# Row 10
$ws.Range("E10").Value = '  +14.49%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.08%  '

# Row 12
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +17.59%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.93%  '

# Row 14
$ws.Range("D14").Value = '2.601.18'
$ws.Range("E14").Value = '  +2.70%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.80%  '

# Row 16
$ws.Range("E16").Value = '  +6.54%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.824'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.83%  '

# Row 18
$ws.Range("D18").Value = '2.270.07'
$ws.Range("E18").Value = '  +2.89%  '

# Row 19
$ws.Range("D19").Value = '43.865.14'
$ws.Range("E19").Value = '  +5.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +13.75%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.37%  '

# Row 22
$ws.Range("E22").Value = '  +0.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.32%  '

# Row 24
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("E25").Value = '  +4.82%  '

# Row 26
$ws.Range("E26").Value = '  +1.58%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.84%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.92%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.80%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.137'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.48%  '

# Row 33
$ws.Range("E33").Value = '  +2.54%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0681'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.91%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.12%  '

# Row 39
$ws.Range("E39").Value = '  +0.25%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0251'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.36%  '

# Row 41
$ws.Range("E41").Value = '  +0.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.17%  '

# Row 43
$ws.Range("E43").Value = '  +8.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0962'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.73%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '97.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.02%  '

# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.32%  '

# Row 47
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.35%  '

# Row 48
$ws.Range("B48").Value = 'TerraClassic'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000209'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -14.27%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.460.79'
$ws.Range("E49").Value = '  -0.08%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.23%  '

# Row 51
$ws.Range("E51").Value = '  +0.52%  '
